$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 704-706; existing rows 704-718 shift down to 707-721
$ws.Range("A704:R706").Insert()

# Row 704
$ws.Range("A704").Value = 6
$ws.Range("B704").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C704").Value = "Metropolitana"
$ws.Range("D704").Value = 44595
$ws.Range("E704").Value = 13
$ws.Range("F704").Value = 100112031
$ws.Range("G704").Value = "Poroto verde"
$ws.Range("H704").Value = "Magnum"
$ws.Range("I704").Value = "Primera"
$ws.Range("J704").Value = 400
$ws.Range("K704").Value = 27000
$ws.Range("L704").Value = 30000
$ws.Range("M704").Value = 28275
$ws.Range("N704").Value = "`$/saco 25 kilos"
$ws.Range("O704").Value = "Región Metropolitana"
$ws.Range("P704").Value = 1131
$ws.Range("Q704").Value = 25
$ws.Range("R704").Value = "Hortaliza"

# Row 705
$ws.Range("A705").Value = 6
$ws.Range("B705").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C705").Value = "Metropolitana"
$ws.Range("D705").Value = 44595
$ws.Range("E705").Value = 13
$ws.Range("F705").Value = 100112031
$ws.Range("G705").Value = "Poroto verde"
$ws.Range("H705").Value = "Magnum"
$ws.Range("I705").Value = "Primera"
$ws.Range("J705").Value = 230
$ws.Range("K705").Value = 27000
$ws.Range("L705").Value = 30000
$ws.Range("M705").Value = 28304
$ws.Range("N705").Value = "`$/saco 25 kilos"
$ws.Range("O705").Value = "Región del Maule"
$ws.Range("P705").Value = 1132
$ws.Range("Q705").Value = 25
$ws.Range("R705").Value = "Hortaliza"

# Row 706
$ws.Range("A706").Value = 6
$ws.Range("B706").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C706").Value = "Metropolitana"
$ws.Range("D706").Value = 44595
$ws.Range("E706").Value = 13
$ws.Range("F706").Value = 100112031
$ws.Range("G706").Value = "Poroto verde"
$ws.Range("H706").Value = "Sin especificar"
$ws.Range("I706").Value = "Primera"
$ws.Range("J706").Value = 180
$ws.Range("K706").Value = 40000
$ws.Range("L706").Value = 45000
$ws.Range("M706").Value = 42222
$ws.Range("N706").Value = "`$/malla 25 kilos"
$ws.Range("O706").Value = "Provincia del Elquí"
$ws.Range("P706").Value = 1689
$ws.Range("Q706").Value = 25
$ws.Range("R706").Value = "Hortaliza"

